# Actualización automática: 2025-07-17 08:10
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

# --- Update row 2 with the new event data ---
$ws.Range("A2").Value = "281474991395097-1752681434684"
$ws.Range("B2").Value = "Harsh Brake"
$ws.Range("C2").Value = "2025-07-16T09:57:14.684"

# Columns D, E, F hold numeric-looking identifiers that must stay as text,
# matching the original sheet's inline-string cell type.
$ws.Range("D2:F2").NumberFormat = "@"
$ws.Range("D2").Value = "281474991395097"
$ws.Range("E2").Value = "125"
$ws.Range("F2").Value = "51834055"
# Drop the text number-format so the cells don't pick up an extra style,
# keeping them stored as plain text like the rest of the sheet.
$ws.Range("D2:F2").Style = "Normal"

$ws.Range("G2").Value = "DAVID SERRANO"
$ws.Range("H2").Value = 20.67662109
$ws.Range("I2").Value = -103.429516539
$ws.Range("J2").Value = 0.7816076278686523
$ws.Range("K2").Value = "No video URL"
$ws.Range("L2").Value = "No video URL"

# --- Remove the two obsolete rows (3 and 4) ---
$ws.Rows("3:4").Delete()
